$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.061.32"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "3.029.66"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'" + "554.95"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").Value = "'" + "136.20"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.020.94"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "'" + "0.498"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -4.43%  "
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("D12").Value = "'" + "0.451"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "'" + "34.40"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "3.516.21"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "62.155.78"
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "3.027.67"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "'" + "6.69"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'" + "474.40"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").Value = "'" + "13.31"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "'" + "0.677"
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").Value = "'" + "7.10"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'" + "80.52"
$ws.Range("D25").Value = "'" + "12.17"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'" + "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'" + "2.74"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").Value = "'" + "7.78"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").Value = "'" + "0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "'" + "25.83"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "'" + "1.16"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "'" + "2.33"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "'" + "55.63"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "'" + "461.17"
$ws.Range("E37").Value = "  -7.40%  "
$ws.Range("D38").Value = "3.226.92"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").Value = "'" + "0.0799"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'" + "0.0386"
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'" + "2.48"
$ws.Range("E43").Value = "  -7.43%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'" + "25.91"
$ws.Range("E45").Value = "  +5.14%  "
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "'" + "0.109"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'" + "118.54"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("D50").Value = "'" + "0.0₃0498"
$ws.Range("E50").Value = "  -6.44%  "
$ws.Range("E51").Value = "  +7.16%  "
